$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.071.36"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.400.90"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "2.408.35"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "60.833.96"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "2.412.79"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "587.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  +6.62%  "
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.75%  "
